# Add the three new mailing-list email addresses that were appended
# to the bottom of the (single-column) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A114").Value = "juliette.soulier@curie.fr"
$ws.Range("A115").Value = "hugo.laporte@curie.fr"
$ws.Range("A116").Value = "maxime.dubail@curie.fr"

# Restore the view/selection state recorded the next time the file was saved.
$excel.ActiveWindow.ScrollRow = 89
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G102").Select()
